# Auto-generated Excel COM-interop script implementing the "Update countries & provincias Spain" diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string (country name) rank swaps -------------------------------
# Colombia overtakes Banglades (row 19/20)
$ws.Range("A19").Value = "Colombia"
$ws.Range("A20").Value = "Banglades"

# Bulgaria overtakes Macedonia/Palestina (rows 80/81/82)
$ws.Range("A80").Value = "Bulgaria"
$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("A82").Value = "Estado de Palestina"

# Islas Malvinas overtakes Groenlandia (rows 210/211)
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"

# --- Updated daily figures --------------------------------------------------
# Row 4
$ws.Range("B4").Value = 4093329
$ws.Range("C4").Value = 64760
$ws.Range("D4").Value = 1931676
$ws.Range("E4").Value = 2015613
$ws.Range("G4").Value = 1087
$ws.Range("H4").Value = 146040

# Row 5
$ws.Range("D5").Value = 1532138
$ws.Range("E5").Value = 612605

# Row 19
$ws.Range("B19").Value = 218428
$ws.Range("C19").Value = 7390
$ws.Range("D19").Value = 101613
$ws.Range("E19").Value = 109442
$ws.Range("G19").Value = 207
$ws.Range("H19").Value = 7373

# Row 20
$ws.Range("B20").Value = 213254
$ws.Range("C20").Value = 2744
$ws.Range("D20").Value = 117202
$ws.Range("E20").Value = 93301
$ws.Range("G20").Value = 42
$ws.Range("H20").Value = 2751

# Row 21
$ws.Range("B21").Value = 204470
$ws.Range("C21").Value = 580
$ws.Range("E21").Value = 6688

# Row 24
$ws.Range("B24").Value = 112206
$ws.Range("C24").Value = 509
$ws.Range("D24").Value = 98127
$ws.Range("E24").Value = 5209
$ws.Range("G24").Value = 8
$ws.Range("H24").Value = 8870

# Row 28
$ws.Range("B28").Value = 89745
$ws.Range("C28").Value = 667
$ws.Range("D28").Value = 30075
$ws.Range("E28").Value = 55230
$ws.Range("G28").Value = 41
$ws.Range("H28").Value = 4440

# Row 48
$ws.Range("B48").Value = 41135
$ws.Range("C48").Value = 906
$ws.Range("D48").Value = 27756
$ws.Range("E48").Value = 11806
$ws.Range("G48").Value = 42
$ws.Range("H48").Value = 1573

# Row 50
$ws.Range("B50").Value = 38344
$ws.Range("C50").Value = 543
$ws.Range("D50").Value = 15815
$ws.Range("E50").Value = 21716
$ws.Range("G50").Value = 8
$ws.Range("H50").Value = 813

# Row 80
$ws.Range("B80").Value = 9584
$ws.Range("C80").Value = 330
$ws.Range("D80").Value = 4643
$ws.Range("E80").Value = 4620
$ws.Range("G80").Value = 8
$ws.Range("H80").Value = 321

# Row 81
$ws.Range("B81").Value = 9547
$ws.Range("C81").Value = 135
$ws.Range("D81").Value = 5071
$ws.Range("E81").Value = 4034
$ws.Range("G81").Value = 10
$ws.Range("H81").Value = 442

# Row 82
$ws.Range("B82").Value = 9398
$ws.Range("C82").Value = 170
$ws.Range("D82").Value = 1950
$ws.Range("E82").Value = 7382
$ws.Range("G82").Value = 2
$ws.Range("H82").Value = 66

# Row 179
$ws.Range("D179").Value = 139
$ws.Range("E179").Value = 5

# --- Footer timestamp ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 00:32"
